$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E4").Value = 33
$ws.Range("F4").Value = 22
$ws.Range("H4").Value = 22
$ws.Range("E5").Value = 97
$ws.Range("F5").Value = 62
$ws.Range("H5").Value = 62
$ws.Range("F6").Value = 22
$ws.Range("H6").Value = 22
$ws.Range("E10").Value = 329
$ws.Range("F10").Value = 158
$ws.Range("H10").Value = 158
$ws.Range("E11").Value = 223
$ws.Range("F11").Value = 127
$ws.Range("H11").Value = 127
$ws.Range("E12").Value = 327
$ws.Range("F12").Value = 188
$ws.Range("H12").Value = 188
$ws.Range("E13").Value = 98
$ws.Range("F13").Value = 52
$ws.Range("H13").Value = 52
$ws.Range("E14").Value = 89
$ws.Range("F14").Value = 44
$ws.Range("H14").Value = 44
$ws.Range("F15").Value = 46
$ws.Range("H15").Value = 46
$ws.Range("E16").Value = 141
$ws.Range("F16").Value = 75
$ws.Range("H16").Value = 75
$ws.Range("E17").Value = 61
$ws.Range("E21").Value = 108
$ws.Range("F21").Value = 57
$ws.Range("H21").Value = 57
$ws.Range("E22").Value = 123
$ws.Range("E23").Value = 146
$ws.Range("F23").Value = 67
$ws.Range("H23").Value = 67
$ws.Range("E24").Value = 144
$ws.Range("F24").Value = 74
$ws.Range("H24").Value = 74
$ws.Range("E25").Value = 165
$ws.Range("F25").Value = 81
$ws.Range("H25").Value = 81
$ws.Range("F26").Value = 56
$ws.Range("H26").Value = 56
$ws.Range("E27").Value = 217
$ws.Range("F27").Value = 114
$ws.Range("H27").Value = 114
$ws.Range("E28").Value = 131
$ws.Range("E29").Value = 130
$ws.Range("F29").Value = 76
$ws.Range("H29").Value = 76
$ws.Range("E30").Value = 150
$ws.Range("F30").Value = 85
$ws.Range("H30").Value = 85
$ws.Range("E32").Value = 132
$ws.Range("F32").Value = 72
$ws.Range("H32").Value = 72
$ws.Range("E33").Value = 208
$ws.Range("F33").Value = 106
$ws.Range("H33").Value = 106
$ws.Range("E34").Value = 152
$ws.Range("F34").Value = 94
$ws.Range("H34").Value = 94
$ws.Range("F35").Value = 62
$ws.Range("H35").Value = 62
$ws.Range("E36").Value = 44
$ws.Range("F36").Value = 30
$ws.Range("H36").Value = 30
$ws.Range("E37").Value = 112
$ws.Range("E38").Value = 67
$ws.Range("E39").Value = 142
$ws.Range("F39").Value = 67
$ws.Range("H39").Value = 67
$ws.Range("E40").Value = 186
$ws.Range("E41").Value = 273
$ws.Range("F41").Value = 123
$ws.Range("H41").Value = 123
$ws.Range("E42").Value = 246
$ws.Range("F42").Value = 132
$ws.Range("H42").Value = 132
$ws.Range("E43").Value = 80
$ws.Range("F43").Value = 41
$ws.Range("H43").Value = 41
$ws.Range("E44").Value = 219
$ws.Range("F44").Value = 114
$ws.Range("H44").Value = 114
$ws.Range("E45").Value = 91
$ws.Range("E46").Value = 211
$ws.Range("F46").Value = 115
$ws.Range("H46").Value = 115
$ws.Range("E47").Value = 308
$ws.Range("F47").Value = 155
$ws.Range("H47").Value = 155
$ws.Range("F48").Value = 59
$ws.Range("H48").Value = 59
$ws.Range("E49").Value = 189
$ws.Range("F49").Value = 87
$ws.Range("H49").Value = 87
$ws.Range("E50").Value = 164
$ws.Range("F50").Value = 63
$ws.Range("H50").Value = 63
$ws.Range("E51").Value = 154
$ws.Range("F52").Value = 7
$ws.Range("H52").Value = 7
